$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2: "S" -> "s"
$ws.Range("C2").Value = "s"

# D2: 4 -> 0
$ws.Range("D2").Value = 0

# H2: 64.34999999999999 -> 0
$ws.Range("H2").Value = 0

# B3: empty -> "s"
$ws.Range("B3").Value = "s"

# C3: empty -> "s"
$ws.Range("C3").Value = "s"
